$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 4 - "Updates Since IETF-104" : re-order three bullet lines and
# split the "ECMP" bullet into two runs ("Addtional" + rest).
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$tr4 = $shp4.TextFrame.TextRange

# Paragraph 3 was "Added loopback measurement mode" -> "Added block number TLV"
$tr4.Paragraphs(3).Runs(1).Text = "Added block number TLV"

# Paragraph 4 was "Added Path Segment ID in Figure 4" -> "Added loopback measurement mode"
$tr4.Paragraphs(4).Runs(1).Text = "Added loopback measurement mode"

# Paragraph 5 was "Added block number TLV" -> "Show Path Segment ID in Figure 4"
$tr4.Paragraphs(5).Runs(1).Text = "Show Path Segment ID in Figure 4"

# Paragraph 7 ("Added details for handling ECMP for SR Policy") is split
# into two runs: "Addtional" + " details for handling ECMP for SR Policy"
$para7 = $tr4.Paragraphs(7)
$firstRun = $para7.Runs(1)
$firstRun.InsertBefore("Addtional") | Out-Null
$para7.Runs(2).Text = " details for handling ECMP for SR Policy"

# ---------------------------------------------------------------------
# Slide 5 - "History of the Draft" : shrink the content placeholder,
# reduce all font sizes from 18pt to 16pt and tweak the last bullet's
# wording.
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(2)
$tr5 = $shp5.TextFrame.TextRange

# "WG document - " -> "WG document, perhaps - "
$tr5.Paragraphs(10).Runs(1).Text = "WG document, perhaps - "

# Every paragraph with text drops from sz=1800 to sz=1600
for ($i = 1; $i -le 10; $i++) {
    $tr5.Paragraphs($i).Font.Size = 16
}

# Shrink the placeholder's height (8229600 x 3288807 EMU == 648 x 258.9612pt)
$shp5.Height = 3288807 / 12700
